$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report header/description text (new run date + event count).
$ws.Range("A1").Value = "Description unknown, completed 06/22/2023 11:07:25 EDT, by WPJTOWN1.The search returned: 7 events."

# Two trace events (HRTX541048 "Placed Actual" and CRDX15008 "Placed Actual") are no
# longer part of the result set, so remove those two rows. After the first delete, the
# old row 9 (HRTX541043) shifts up to row 8, so deleting row 7 again removes the old
# CRDX15008 row that shifted up to row 7.
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(7).Delete()

# The remaining HRTX541043 row (now row 7) now reflects a newer "Departure" event out of
# HUTCHINSON, KS (instead of the old "Arrive In-Transit" out of KANSAS CITY, KS), so move
# it up above CRDX15803 and refresh its event details (weights are unchanged).
$ws.Range("A6").Value = "HRTX"
$ws.Range("B6").Value = 541043
$ws.Range("C6").Value = "HUTCHINSON"
$ws.Range("D6").Value = "KS"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 22
$ws.Range("G6").Value = 845
$ws.Range("H6").Value = "Departure"
$ws.Range("I6").Value = "HKCKDE"
$ws.Range("J6").Value = "LOVELAND"
$ws.Range("K6").Value = "CO"
$ws.Range("L6").Value = 258850
$ws.Range("M6").Value = 64200
$ws.Range("N6").Value = 194650
$ws.Range("O6").Value = "HRTX541043"

$ws.Range("A7").Value = "CRDX"
$ws.Range("B7").Value = 15803
$ws.Range("C7").Value = "JOHNSTOWN"
$ws.Range("D7").Value = "CO"
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 15
$ws.Range("G7").Value = 1435
$ws.Range("H7").Value = "Placed Actual"
$ws.Range("I7").ClearContents()
$ws.Range("J7").Value = "LOVELAND"
$ws.Range("K7").Value = "CO"
$ws.Range("L7").Value = 284700
$ws.Range("M7").Value = 66900
$ws.Range("N7").Value = 217800
$ws.Range("O7").Value = "CRDX15803"

# Re-select the Car_no column now that the result set only spans down to row 9.
[void]$ws.Range("O3:O9").Select()
